$wb = $excel.ActiveWorkbook

# ALC row 2 (G=5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 442.5
$ws.Range("I2").Value = 131.66667
$ws.Range("K2").Value = 131.66667
$ws.Range("M2").Value = -18.66667000000001

# ALC row 18 (G=5471)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1298.6666
$ws.Range("I18").Value = 1298.6666
$ws.Range("K18").Value = 1298.6666
$ws.Range("M18").Value = -1014.6666

# ALC row 28 (G=27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 703.6
$ws.Range("I28").Value = 108
$ws.Range("J28").Value = 1597
$ws.Range("K28").Value = 108
$ws.Range("L28").Value = 1597
$ws.Range("M28").Value = 377
$ws.Range("N28").Value = -2567

# ALC row 40 (G=5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5046.1333
$ws.Range("I40").Value = 3799.4443
$ws.Range("J40").Value = 6916.1665
$ws.Range("K40").Value = 3799.4443
$ws.Range("L40").Value = 6916.1665
$ws.Range("M40").Value = -3624.4443
$ws.Range("N40").Value = -7266.1665

# ALC row 41 (G=5478)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2259.8333
$ws.Range("J41").Value = 823.3333
$ws.Range("L41").Value = 823.3333
$ws.Range("N41").Value = -1703.3333

# ALC row 53 (G=5479)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1012.2857
$ws.Range("I53").Value = 138.22223
$ws.Range("K53").Value = 138.22223
$ws.Range("M53").Value = 498.77777

# ALC row 106 (G=19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2410.7646
$ws.Range("I106").Value = 2544
$ws.Range("J106").Value = 2166.5
$ws.Range("K106").Value = 2544
$ws.Range("L106").Value = 2166.5
$ws.Range("M106").Value = -1913
$ws.Range("N106").Value = -3428.5

# ALC row 111 (G=27768)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2897.077
$ws.Range("I111").Value = 3236.8333
$ws.Range("J111").Value = 2605.8572
$ws.Range("K111").Value = 9710.499899999999
$ws.Range("L111").Value = 7817.571599999999
$ws.Range("M111").Value = -6643.499899999999
$ws.Range("N111").Value = -13951.5716

# ALC row 138 (G=44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1982.5
$ws.Range("J138").Value = 3282.4443
$ws.Range("L138").Value = 9847.332900000001
$ws.Range("N138").Value = -20127.3329

# ARM row 32 (G=44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2686.0488
$ws.Range("I32").Value = 1973.0555
$ws.Range("K32").Value = 1973.0555
$ws.Range("M32").Value = -1686.0555

# ARM row 88 (G=12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1848.9565
$ws.Range("I88").Value = 1519.6666
$ws.Range("J88").Value = 2060.6428
$ws.Range("K88").Value = 1519.6666
$ws.Range("L88").Value = 2060.6428
$ws.Range("M88").Value = -1113.6666
$ws.Range("N88").Value = -2872.6428

# ARM row 91 (G=12530)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1848.9565
$ws.Range("I91").Value = 1519.6666
$ws.Range("J91").Value = 2060.6428
$ws.Range("K91").Value = 1519.6666
$ws.Range("L91").Value = 2060.6428
$ws.Range("M91").Value = -115.6666
$ws.Range("N91").Value = -4868.6428

# ARM row 122 (G=36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7928.2
$ws.Range("I122").Value = 7183.2856
$ws.Range("K122").Value = 21549.8568
$ws.Range("M122").Value = -19099.8568

# BSM row 22 (G=5092)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1488.8
$ws.Range("I22").Value = 1661.875
$ws.Range("J22").Value = 1181.1111
$ws.Range("K22").Value = 1661.875
$ws.Range("L22").Value = 1181.1111
$ws.Range("M22").Value = -1488.875
$ws.Range("N22").Value = -1527.1111

# BSM row 86 (G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3769.5
$ws.Range("I86").Value = 3783.875
$ws.Range("K86").Value = 3783.875
$ws.Range("M86").Value = -2660.875

# BSM row 89 (G=12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3769.5
$ws.Range("I89").Value = 3783.875
$ws.Range("K89").Value = 18919.375
$ws.Range("M89").Value = -13303.375

# CRP row 22 (G=5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5964.8945
$ws.Range("I22").Value = 7309.2666
$ws.Range("K22").Value = 7309.2666
$ws.Range("M22").Value = -6959.2666

# CRP row 31 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4738.24
$ws.Range("I31").Value = 3227.8
$ws.Range("J31").Value = 10780
$ws.Range("K31").Value = 3227.8
$ws.Range("L31").Value = 10780
$ws.Range("M31").Value = -2932.8
$ws.Range("N31").Value = -11370

# CRP row 34 (G=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4738.24
$ws.Range("I34").Value = 3227.8
$ws.Range("J34").Value = 10780
$ws.Range("K34").Value = 3227.8
$ws.Range("L34").Value = 10780
$ws.Range("M34").Value = -3025.8
$ws.Range("N34").Value = -11184

# CRP row 51 (G=2039)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 45000
$ws.Range("J51").Value = 45000
$ws.Range("L51").Value = 45000
$ws.Range("N51").Value = -46472

# CRP row 58 (G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22735278
$ws.Range("I58").Value = 29421270
$ws.Range("K58").Value = 29421270
$ws.Range("M58").Value = -29421067

# CRP row 61 (G=2039)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 45000
$ws.Range("J61").Value = 45000
$ws.Range("L61").Value = 45000
$ws.Range("N61").Value = -45696

# CRP row 64 (G=10610)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 46475
$ws.Range("I64").Value = 85900
$ws.Range("J64").Value = 33333.332
$ws.Range("K64").Value = 85900
$ws.Range("L64").Value = 33333.332
$ws.Range("M64").Value = -85652
$ws.Range("N64").Value = -33829.332

# CRP row 67 (G=10610)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 46475
$ws.Range("I67").Value = 85900
$ws.Range("J67").Value = 33333.332
$ws.Range("K67").Value = 85900
$ws.Range("L67").Value = 33333.332
$ws.Range("M67").Value = -85042
$ws.Range("N67").Value = -35049.332

# CRP row 80 (G=12015)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 59999
$ws.Range("J80").Value = 59999
$ws.Range("L80").Value = 59999
$ws.Range("N80").Value = -62245

# CRP row 83 (G=12015)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 59999
$ws.Range("J83").Value = 59999
$ws.Range("L83").Value = 179997
$ws.Range("N83").Value = -191229

# CRP row 105 (G=19928)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1298.25
$ws.Range("I105").Value = 1341
$ws.Range("K105").Value = 1341
$ws.Range("M105").Value = 406

# CRP row 122 (G=36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1826.6875
$ws.Range("I122").Value = 2019.75
$ws.Range("J122").Value = 1247.5
$ws.Range("K122").Value = 6059.25
$ws.Range("L122").Value = 3742.5
$ws.Range("M122").Value = -3609.25
$ws.Range("N122").Value = -8642.5

# CRP row 132 (G=44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 38464720
$ws.Range("I132").Value = 52634410
$ws.Range("J132").Value = 4145.857
$ws.Range("K132").Value = 157903230
$ws.Range("L132").Value = 12437.571
$ws.Range("M132").Value = -157900700
$ws.Range("N132").Value = -17497.571

# CRP row 134 (G=44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 16669698
$ws.Range("I134").Value = 25002426
$ws.Range("K134").Value = 75007278
$ws.Range("M134").Value = -75004743

# CRP row 136 (G=44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 22735278
$ws.Range("I136").Value = 29421270
$ws.Range("K136").Value = 88263810
$ws.Range("M136").Value = -88261260

# CUL row 93 (G=19808)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 15333.167
$ws.Range("J93").Value = 15333.167
$ws.Range("L93").Value = 45999.501
$ws.Range("N93").Value = -49743.501

# CUL row 131 (G=36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2030.7142
$ws.Range("I131").Value = 2057.25
$ws.Range("J131").Value = 1995.3334
$ws.Range("K131").Value = 6171.75
$ws.Range("L131").Value = 5986.0002
$ws.Range("M131").Value = -1131.75
$ws.Range("N131").Value = -16066.0002

# GSM row 70 (G=14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5976.8237
$ws.Range("I70").Value = 5764.857
$ws.Range("K70").Value = 5764.857
$ws.Range("M70").Value = -5494.857

# GSM row 73 (G=14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5976.8237
$ws.Range("I73").Value = 5764.857
$ws.Range("K73").Value = 5764.857
$ws.Range("M73").Value = -4828.857

# LTW row 4 (G=3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 11666.667
$ws.Range("J4").Value = 11666.667
$ws.Range("L4").Value = 11666.667
$ws.Range("N4").Value = -11892.667

# LTW row 5 (G=3790)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 20000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -19887
$ws.Range("N5").Value = -10226

# LTW row 22 (G=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2557.9092
$ws.Range("I22").Value = 2554.75
$ws.Range("J22").Value = 2566.3333
$ws.Range("K22").Value = 2554.75
$ws.Range("L22").Value = 2566.3333
$ws.Range("M22").Value = -2259.75
$ws.Range("N22").Value = -3156.3333

# LTW row 27 (G=5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2557.9092
$ws.Range("I27").Value = 2554.75
$ws.Range("J27").Value = 2566.3333
$ws.Range("K27").Value = 2554.75
$ws.Range("L27").Value = 2566.3333
$ws.Range("M27").Value = -2447.75
$ws.Range("N27").Value = -2780.3333

# LTW row 28 (G=3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 11666.667
$ws.Range("J28").Value = 11666.667
$ws.Range("L28").Value = 11666.667
$ws.Range("N28").Value = -12130.667

# LTW row 37 (G=3788)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 11666.667
$ws.Range("J37").Value = 11666.667
$ws.Range("L37").Value = 11666.667
$ws.Range("N37").Value = -11880.667

# LTW row 46 (G=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2451.4443
$ws.Range("I46").Value = 2451.4443
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2451.4443
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2263.4443

# LTW row 55 (G=5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 499.47827
$ws.Range("I55").Value = 307.8889
$ws.Range("J55").Value = 622.6429000000001
$ws.Range("K55").Value = 307.8889
$ws.Range("L55").Value = 622.6429000000001
$ws.Range("M55").Value = -134.8889
$ws.Range("N55").Value = -968.6429000000001

# WVR row 87 (G=12005)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 49999
$ws.Range("J87").Value = 49999
$ws.Range("L87").Value = 49999
$ws.Range("N87").Value = -52495

# WVR row 90 (G=12005)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 49999
$ws.Range("J90").Value = 49999
$ws.Range("L90").Value = 149997
$ws.Range("N90").Value = -162477
